$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
[void]$tbl.Rows.Add()
$tbl = $d.Tables.Item(1)
$newRowIndex = $tbl.Rows.Count
$tbl.Cell($newRowIndex, 1).Range.Text = "2.0"
$tbl.Cell($newRowIndex, 2).Range.Text = "23/05/2016"
$tbl.Cell($newRowIndex, 3).Range.Text = "Versão Final"
$tbl.Cell($newRowIndex, 4).Range.Text = "Rogério"

$cellRange = $tbl.Cell($newRowIndex, 4).Range
$bookmarkPos = $cellRange.End - 1
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
[void]$d.Bookmarks.Add("TestBM", $bmRange)
$check = $d.Bookmarks.Item("TestBM").Range
Write-Output ("TestBM range Start/End: " + $check.Start + "/" + $check.End)
